# Fill in invoice form data (test/placeholder text entered via the UI)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---
$ws.Range("I5").Value = "gj"
$ws.Range("K5").Value = "Date: 02-Mar-22 - (Wednesday)"

$ws.Range("B6").Value = "g"
$ws.Range("G6").Value = "  gfhgg"
$ws.Range("J6").Value = "gjh"
$ws.Range("N6").Value = "gj"

# --- Goods table (rows 10-12 filled in, row 13-17 left blank) ---
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "hg"
$ws.Range("G10").Value = "jhg"
$ws.Range("H10").Value = "jh"
$ws.Range("I10").Value = "g"
$ws.Range("K10").Value = "jhg"
$ws.Range("M10").Value = "gjh"

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "gjh"
$ws.Range("G11").Value = "hg"
$ws.Range("H11").Value = "jhg"
$ws.Range("I11").Value = "jhg"
$ws.Range("K11").Value = "jh"

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "g"

# --- Old jewellery table (rows 20-22) ---
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "h"
$ws.Range("I20").Value = "hjh"
$ws.Range("J20").Value = "jh"
$ws.Range("M20").Value = "gjh"

$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "g"
$ws.Range("I21").Value = "jhg"
$ws.Range("J21").Value = "jh"
$ws.Range("M21").Value = "gjh"

$ws.Range("A22").Value = 3
$ws.Range("B22").Value = "gj"
$ws.Range("I22").Value = "hg"
$ws.Range("J22").Value = "jhg"
$ws.Range("M22").Value = "jhg"

# --- Other Addition Or Deduction table (rows 24-27) ---
$ws.Range("D24").Value = "Other Addition Or Deduction"
$ws.Range("M24").Value = "Amount"

$ws.Range("A25").Value = 1
$ws.Range("D25").Value = "bhkj"
$ws.Range("M25").Value = "hj"

$ws.Range("A26").Value = 2
$ws.Range("D26").Value = "jk"
$ws.Range("M26").Value = "hkj"

$ws.Range("A27").Value = 3
$ws.Range("D27").Value = "hkj"
$ws.Range("M27").Value = "h"
